$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-24 20:47:19"

$wsZhCn.Range("H4").Value = "2016-08-24 20:47:14"
$wsZhCn.Range("K4").Value = "2016-08-24 20:47:41"

$wsDeDe.Range("H4").Value = "2016-08-24 20:47:19"
$wsDeDe.Range("K4").Value = "2016-08-24 20:47:49"
